$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5483334058394633
$ws.Range("B3").Value = 0.4388152899725493
$ws.Range("B4").Value = 0.2924055799924064
$ws.Range("B5").Value = 0.2633722774678172
$ws.Range("B6").Value = 0.7261509110052229
$ws.Range("B7").Value = 0.5856228053907445
$ws.Range("B8").Value = 0.5279640510362523
$ws.Range("B9").Value = 0.4028384916812343
$ws.Range("B10").Value = 0.3070763867582285
$ws.Range("B11").Value = 0.04688115307313914
$ws.Range("B12").Value = 0.6404433048987681
$ws.Range("B13").Value = 0.5512165644000019
$ws.Range("B14").Value = 0.527506874197718
$ws.Range("B15").Value = 0.3362693909433243
